# Generate Report for Archive
# - Status cells move from "Ready for handoff" to "In Translation"
#   (Overview!E2:F2, zh-cn!C2, de-de!C2)
# - The now-narrower "Status" column is resized on each sheet that shows it

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"
$newWidth = 13.4101845877511

# Overview sheet: Status text lives in columns E (zh-cn) and F (de-de)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = $newWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth

# zh-cn sheet: Status text lives in column C
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Columns.Item(3).ColumnWidth = $newWidth

# de-de sheet: Status text lives in column C
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Columns.Item(3).ColumnWidth = $newWidth
